# Update the "datetimeFigureOut" date field text from 2024/11/20 to 2024/12/3
# across the slide master and every slide layout's date placeholder.

$p = $ppt.ActivePresentation

# Slide master date placeholder (Shapes.Item(3))
$m = $p.SlideMaster
$m.Shapes.Item(3).TextFrame.TextRange.Text = "2024/12/3"

# Slide layouts: map of layout index -> date-placeholder shape index
$layoutDateShapeIndex = @{
    1 = 3
    2 = 3
    3 = 3
    4 = 4
    5 = 6
    6 = 2
    7 = 1
    8 = 4
    9 = 4
    10 = 3
    11 = 3
}

$layouts = $m.CustomLayouts
foreach ($li in $layoutDateShapeIndex.Keys) {
    $shpIdx = $layoutDateShapeIndex[$li]
    $layout = $layouts.Item($li)
    $layout.Shapes.Item($shpIdx).TextFrame.TextRange.Text = "2024/12/3"
}

# Slide 1: fix the "s1" rounded-rectangle shape's run/endParaRPr properties
$s = $p.Slides.Item(1)
$s1Shape = $s.Shapes.Item(2)
$s1Shape.TextFrame.TextRange.Text = "s1"
